{"js": "// Replace the date string and every \"a\u00f7b=c, d\" practice-problem answer in\n// the document, in document order. A couple of the \"before\" values are not\n// unique (e.g. \"90\u00f74=22, 2\" shows up twice), so each pair is applied by\n// searching for the (still unmodified) old text and rewriting only the\n// first match \u2014 once that occurrence is rewritten it no longer matches,\n// so the next lookup of the same old text (if any) naturally lands on the\n// next remaining occurrence.\nconst replacements = [\n  [\"2024-07-22 Monday\", \"2024-07-23 Tuesday\"],\n  [\"33\u00f76=5, 3\", \"73\u00f74=18, 1\"],\n  [\"94\u00f78=11, 6\", \"62\u00f77=8, 6\"],\n  [\"40\u00f76=6, 4\", \"33\u00f73=11, 0\"],\n  [\"88\u00f74=22, 0\", \"33\u00f73=11, 0\"],\n  [\"90\u00f79=10, 0\", \"30\u00f72=15, 0\"],\n  [\"89\u00f78=11, 1\", \"13\u00f78=1, 5\"],\n  [\"54\u00f79=6, 0\", \"29\u00f73=9, 2\"],\n  [\"10\u00f75=2, 0\", \"37\u00f76=6, 1\"],\n  [\"90\u00f74=22, 2\", \"60\u00f75=12, 0\"],\n  [\"69\u00f76=11, 3\", \"72\u00f73=24, 0\"],\n  [\"14\u00f76=2, 2\", \"27\u00f74=6, 3\"],\n  [\"68\u00f75=13, 3\", \"33\u00f74=8, 1\"],\n  [\"62\u00f76=10, 2\", \"62\u00f78=7, 6\"],\n  [\"79\u00f72=39, 1\", \"14\u00f79=1, 5\"],\n  [\"17\u00f72=8, 1\", \"85\u00f74=21, 1\"],\n  [\"90\u00f74=22, 2\", \"56\u00f76=9, 2\"],\n  [\"72\u00f78=9, 0\", \"66\u00f78=8, 2\"],\n  [\"58\u00f72=29, 0\", \"79\u00f77=11, 2\"],\n  [\"33\u00f78=4, 1\", \"55\u00f73=18, 1\"],\n  [\"45\u00f77=6, 3\", \"86\u00f74=21, 2\"],\n  [\"81\u00f74=20, 1\", \"65\u00f74=16, 1\"],\n  [\"31\u00f79=3, 4\", \"81\u00f77=11, 4\"],\n  [\"62\u00f74=15, 2\", \"63\u00f72=31, 1\"],\n  [\"97\u00f77=13, 6\", \"78\u00f79=8, 6\"],\n  [\"20\u00f76=3, 2\", \"55\u00f77=7, 6\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${oldText}\"`);\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Replace the date string and every \"a\u00f7b=c, d\" practice-problem answer in\n# the document, in document order. A couple of the \"before\" values are not\n# unique (e.g. \"90\u00f74=22, 2\" shows up twice), so each pair is applied with\n# Find.Execute(..., Replace:=wdReplaceOne) over the FULL document range\n# (re-fetched fresh each iteration so the search restarts at the top).\n# Replacing just the first hit means once that occurrence is rewritten it\n# no longer matches, so a later lookup of the same old text (if any)\n# naturally lands on the next remaining occurrence.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-07-22 Monday\", \"2024-07-23 Tuesday\"),\n    @(\"33\u00f76=5, 3\", \"73\u00f74=18, 1\"),\n    @(\"94\u00f78=11, 6\", \"62\u00f77=8, 6\"),\n    @(\"40\u00f76=6, 4\", \"33\u00f73=11, 0\"),\n    @(\"88\u00f74=22, 0\", \"33\u00f73=11, 0\"),\n    @(\"90\u00f79=10, 0\", \"30\u00f72=15, 0\"),\n    @(\"89\u00f78=11, 1\", \"13\u00f78=1, 5\"),\n    @(\"54\u00f79=6, 0\", \"29\u00f73=9, 2\"),\n    @(\"10\u00f75=2, 0\", \"37\u00f76=6, 1\"),\n    @(\"90\u00f74=22, 2\", \"60\u00f75=12, 0\"),\n    @(\"69\u00f76=11, 3\", \"72\u00f73=24, 0\"),\n    @(\"14\u00f76=2, 2\", \"27\u00f74=6, 3\"),\n    @(\"68\u00f75=13, 3\", \"33\u00f74=8, 1\"),\n    @(\"62\u00f76=10, 2\", \"62\u00f78=7, 6\"),\n    @(\"79\u00f72=39, 1\", \"14\u00f79=1, 5\"),\n    @(\"17\u00f72=8, 1\", \"85\u00f74=21, 1\"),\n    @(\"90\u00f74=22, 2\", \"56\u00f76=9, 2\"),\n    @(\"72\u00f78=9, 0\", \"66\u00f78=8, 2\"),\n    @(\"58\u00f72=29, 0\", \"79\u00f77=11, 2\"),\n    @(\"33\u00f78=4, 1\", \"55\u00f73=18, 1\"),\n    @(\"45\u00f77=6, 3\", \"86\u00f74=21, 2\"),\n    @(\"81\u00f74=20, 1\", \"65\u00f74=16, 1\"),\n    @(\"31\u00f79=3, 4\", \"81\u00f77=11, 4\"),\n    @(\"62\u00f74=15, 2\", \"63\u00f72=31, 1\"),\n    @(\"97\u00f77=13, 6\", \"78\u00f79=8, 6\"),\n    @(\"20\u00f76=3, 2\", \"55\u00f77=7, 6\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 1)\n\n    if (-not $found) {\n        throw \"Could not find text to replace: '$oldText'\"\n    }\n}\n"}
